$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.174.01"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.603.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.68"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.82"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000305"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.177.40"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.10"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "592.71"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.296.75"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.602.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.77%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.74"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.55%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.66"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.02"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.93"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.77"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.66%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.90%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.22"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0895"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.957.10"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.13%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.63%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "522.80"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.89"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.51%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000252"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.86%  "
